$d = $word.ActiveDocument

# 1. Update the document title
$d.Content.Find.Execute("Meetrapport titel", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Hoe ziet een wasprogramma eruit?", 2)

# 2. Re-join the "Een duidelijke opsoming ..." sentence that had been split
#    across runs by a spell-check proofErr marker around "opsoming".
$d.Content.Find.Execute("Een duidelijke opsoming maken van de stappen in een was programma.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Een duidelijke opsoming maken van de stappen in een was programma.", 2)

# 3. Re-join the "Deze website geeft ..." sentence that had been split
#    across runs by a spell-check proofErr marker around "uitgebreiden".
$d.Content.Find.Execute("Deze website geeft een duidelijk en uitgebreiden uitleg over hoe een wasprogramma er uit ziet en wat de verschillen kunnen zijn. Hier staat ook wat er nodig is om dit uit te kunnen voeren.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Deze website geeft een duidelijk en uitgebreiden uitleg over hoe een wasprogramma er uit ziet en wat de verschillen kunnen zijn. Hier staat ook wat er nodig is om dit uit te kunnen voeren.", 2)

# 4. Re-join the "Uit onderzoek is een lijst ..." sentence that had been
#    split across runs by a spell-check proofErr marker around "bonte-was".
$d.Content.Find.Execute("Uit onderzoek is een lijst van stappen gekomen waar de wasmachine door heen moet om een goede was te draaien. Hier zijn ook dingen uit gekomen om rekening bijvoorbeeld het toerental waar de trommel op draait tijdens het centrifugeren (lager voor fijnere was dan voor bonte was), het spoelen van de was verschild ook tussen bonte en fijne was zo moet bonte-was wel 6 keer en fijne was maar 3 keer gespoeld worden.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Uit onderzoek is een lijst van stappen gekomen waar de wasmachine door heen moet om een goede was te draaien. Hier zijn ook dingen uit gekomen om rekening bijvoorbeeld het toerental waar de trommel op draait tijdens het centrifugeren (lager voor fijnere was dan voor bonte was), het spoelen van de was verschild ook tussen bonte en fijne was zo moet bonte-was wel 6 keer en fijne was maar 3 keer gespoeld worden.", 2)

# 5. Re-join the "Spoelen (bont-was ..." text that had been split across
#    runs by a spell-check proofErr marker around "bont-was".
$d.Content.Find.Execute("Spoelen (bont-was wel 6 keer fijne was 3 keer)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Spoelen (bont-was wel 6 keer fijne was 3 keer)", 2)
